# Added date/datetime support to excel parser
#
# Adds two new columns to Sheet1:
#   D = "some_date"      (date-only values, formatted D/M/YY)
#   E = "some_datetime"  (date+time values, formatted DD/MM/YYYY HH:MM:SS)
# and re-formats the existing "some_time" column (C) to HH:MM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers -----------------------------------------------------
$ws.Range("D1").Value() = "some_date"
$ws.Range("E1").Value() = "some_datetime"

# --- Row 2 -------------------------------------------------------------
# some_date: 2010-10-10, some_datetime: 2010-10-10 10:30:10
$ws.Range("D2").Value() = 40461
$ws.Range("D2").NumberFormat() = "D/M/YY"

$ws.Range("E2").Value() = 40461.4376157407
$ws.Range("E2").NumberFormat() = "DD/MM/YYYY\ HH:MM:SS"

# --- Row 3 -------------------------------------------------------------
# some_date: 2010-11-01, some_datetime: 2010-11-01 07:40:00
$ws.Range("D3").Value() = 40483
$ws.Range("D3").NumberFormat() = "D/M/YY"

$ws.Range("E3").Value() = 40483.3194444444
$ws.Range("E3").NumberFormat() = "DD/MM/YYYY\ HH:MM:SS"

# --- Re-format existing time column (C) from "HH:MM:SS AM/PM" to "HH:MM"
$ws.Range("C2:C3").NumberFormat() = "HH:MM"

# --- Column widths (best-effort match to target layout) ----------------
$ws.Columns.Item(1).ColumnWidth() = 13.333333333333334
$ws.Columns.Item(2).ColumnWidth() = 11.5
$ws.Columns.Item(3).ColumnWidth() = 17.0
$ws.Columns.Item(4).ColumnWidth() = 19.0
$ws.Columns.Item(5).ColumnWidth() = 21.833333333333332

# --- Move the active selection to D3, matching the edited workbook -----
$ws.Range("D3").Select() | Out-Null
